$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update CMIP6 Specialization ID mappings (column C) ---
# New shared strings must be introduced in this exact order so that the
# resulting sharedStrings table indices (40 .. 53) line up with the target
# workbook; the order follows the target index sequence, not sheet order.

$ws.Range("C5").Value  = "cmip6.seaice.grid.discretisation.vertical.layering"                     # -> index 40
$ws.Range("C9").Value  = "cmip6.seaice.grid.discretisation.horizontal.thermodynamics_time_step"   # -> index 41
$ws.Range("C7").Value  = "cmip6.seaice.grid.discretisation.horizontal.dynamics_time_step"         # -> index 42
$ws.Range("C16").Value = "cmip6.seaice.grid.discretisation.horizontal.grid"                       # -> index 43
$ws.Range("C14").Value = "cmip6.seaice.grid.discretisation.horizontal.scheme"                     # -> index 44
$ws.Range("C19").Value = "cmip6.seaice.dynamics.redistribution"                                   # -> index 45
$ws.Range("C32").Value = "cmip6.seaice.radiative_processes.surface_albedo"                        # -> index 46
$ws.Range("C33").Value = "cmip6.seaice.thermodynamics.melt_ponds.formulation"                     # -> index 47
$ws.Range("C23").Value = "cmip6.seaice.thermodynamics.energy.basal_heat_flux"                     # -> index 48
$ws.Range("C25").Value = "cmip6.seaice.dynamics.transport_in_thickness_space"                     # -> index 49
$ws.Range("C26").Value = "cmip6.seaice.thermodynamics.energy.heat_diffusion"                      # -> index 50
$ws.Range("C27").Value = "cmip6.seaice.thermodynamics.mass.new_ice_formation"                     # -> index 51
$ws.Range("C29").Value = "cmip6.seaice.thermodynamics.snow_processes.redistribution"              # -> index 52
$ws.Range("C15").Value = "cmip6.seaice.dynamics.horizontal_transport"                             # -> index 53

# Row 31 picks up the string that used to live in row 29 (already index 39,
# reused rather than duplicated).
$ws.Range("C31").Value = "cmip6.seaice.thermodynamics.snow_processes.snow_ice_formation_scheme"

# --- Apply the "Normal" style (clears the row border format) on the cells
#     whose style moves from the bordered style to the new plain style ---
$ws.Range("C5").Style = "Normal"
$ws.Range("C6").Style = "Normal"
$ws.Range("C7").Style = "Normal"
$ws.Range("C14").Style = "Normal"
$ws.Range("C15").Style = "Normal"
$ws.Range("C16").Style = "Normal"
$ws.Range("C18").Style = "Normal"
$ws.Range("C23").Style = "Normal"
$ws.Range("C25").Style = "Normal"
$ws.Range("C26").Style = "Normal"
$ws.Range("C27").Style = "Normal"
$ws.Range("C29").Style = "Normal"
$ws.Range("C30").Style = "Normal"
$ws.Range("C31").Style = "Normal"
$ws.Range("C32").Style = "Normal"
$ws.Range("C33").Style = "Normal"

# --- Update selected cell shown in the sheet view ---
$ws.Range("C13").Select()
